{"js": "// Update the date line and all \"A\u00f7B=\" problems to the new values, per the\n// commit's regenerated worksheet content. Each old value is unique in the\n// document, so a simple text search + replace (first match) per pair is\n// safe and order-independent.\nconst replacements = [\n  [\"2024-04-20 Saturday\", \"2024-04-21 Sunday\"],\n  [\"426\u00f76=\", \"762\u00f77=\"],\n  [\"601\u00f74=\", \"143\u00f76=\"],\n  [\"649\u00f75=\", \"145\u00f77=\"],\n  [\"401\u00f77=\", \"479\u00f75=\"],\n  [\"398\u00f73=\", \"616\u00f78=\"],\n  [\"617\u00f72=\", \"589\u00f73=\"],\n  [\"895\u00f73=\", \"992\u00f76=\"],\n  [\"657\u00f76=\", \"941\u00f77=\"],\n  [\"653\u00f79=\", \"812\u00f77=\"],\n  [\"571\u00f74=\", \"283\u00f75=\"],\n  [\"682\u00f72=\", \"561\u00f75=\"],\n  [\"561\u00f72=\", \"456\u00f73=\"],\n  [\"600\u00f78=\", \"857\u00f76=\"],\n  [\"867\u00f77=\", \"176\u00f73=\"],\n  [\"299\u00f75=\", \"455\u00f73=\"],\n  [\"877\u00f75=\", \"582\u00f78=\"],\n  [\"953\u00f79=\", \"477\u00f74=\"],\n  [\"698\u00f77=\", \"754\u00f75=\"],\n  [\"866\u00f76=\", \"947\u00f79=\"],\n  [\"289\u00f74=\", \"177\u00f78=\"],\n  [\"439\u00f72=\", \"679\u00f72=\"],\n  [\"819\u00f75=\", \"157\u00f75=\"],\n  [\"547\u00f74=\", \"138\u00f79=\"],\n  [\"802\u00f75=\", \"551\u00f72=\"],\n  [\"538\u00f75=\", \"792\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and all \"A\u00f7B=\" problems to the new values, per the\n# commit's regenerated worksheet content. Each old value is unique in the\n# document, so Find/Replace (one match each) per pair is safe and\n# order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-20 Saturday\", \"2024-04-21 Sunday\"),\n    @(\"426\u00f76=\", \"762\u00f77=\"),\n    @(\"601\u00f74=\", \"143\u00f76=\"),\n    @(\"649\u00f75=\", \"145\u00f77=\"),\n    @(\"401\u00f77=\", \"479\u00f75=\"),\n    @(\"398\u00f73=\", \"616\u00f78=\"),\n    @(\"617\u00f72=\", \"589\u00f73=\"),\n    @(\"895\u00f73=\", \"992\u00f76=\"),\n    @(\"657\u00f76=\", \"941\u00f77=\"),\n    @(\"653\u00f79=\", \"812\u00f77=\"),\n    @(\"571\u00f74=\", \"283\u00f75=\"),\n    @(\"682\u00f72=\", \"561\u00f75=\"),\n    @(\"561\u00f72=\", \"456\u00f73=\"),\n    @(\"600\u00f78=\", \"857\u00f76=\"),\n    @(\"867\u00f77=\", \"176\u00f73=\"),\n    @(\"299\u00f75=\", \"455\u00f73=\"),\n    @(\"877\u00f75=\", \"582\u00f78=\"),\n    @(\"953\u00f79=\", \"477\u00f74=\"),\n    @(\"698\u00f77=\", \"754\u00f75=\"),\n    @(\"866\u00f76=\", \"947\u00f79=\"),\n    @(\"289\u00f74=\", \"177\u00f78=\"),\n    @(\"439\u00f72=\", \"679\u00f72=\"),\n    @(\"819\u00f75=\", \"157\u00f75=\"),\n    @(\"547\u00f74=\", \"138\u00f79=\"),\n    @(\"802\u00f75=\", \"551\u00f72=\"),\n    @(\"538\u00f75=\", \"792\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
